$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.599.08"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +4.09%  '
$ws.Range('D3').Value = "'1.743.32"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.47%  '
$ws.Range('D4').Value = "'0.9993"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'246.39"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.84%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').Value = "'0.4827"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.90%  '
$ws.Range('D8').Value = "'0.2693"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.45%  '
$ws.Range('D9').Value = "'0.06262"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.47%  '
$ws.Range('D10').Value = "'1.743.71"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.52%  '
$ws.Range('D11').Value = "'0.07137"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.84%  '
$ws.Range('D12').Value = "'15.91"
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Value = "'0.6251"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.71%  '
$ws.Range('D14').Value = "'4.524"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.77%  '
$ws.Range('D15').Value = "'77.49"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.86%  '
$ws.Range('D16').Value = "'0.9998"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('D17').Value = "'26.604.74"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.14%  '
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').Value = "'0.000006916"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.76%  '
$ws.Range('D20').Value = "'11.79"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.36%  '
$ws.Range('D21').Value = "'1.968.30"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.56%  '
$ws.Range('D22').Value = "'4.625"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.26%  '
$ws.Range('D23').Value = "'8.888"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.26%  '
$ws.Range('D24').Value = "'5.373"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.81%  '
$ws.Range('D25').Value = "'136.05"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.76%  '
$ws.Range('E26').Value = '  +2.64%  '
$ws.Range('D27').Value = "'1.817"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.95%  '
$ws.Range('D28').Value = "'1.433"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.96%  '
$ws.Range('D29').Value = "'106.79"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.28%  '
$ws.Range('D30').Value = "'4.017"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.67%  '
$ws.Range('D31').Value = "'3.744"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.39%  '
$ws.Range('D32').Value = "'0.07892"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.83%  '
$ws.Range('D33').Value = "'0.04595"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.60%  '
$ws.Range('D34').Value = "'2.617"
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Value = "'0.6416"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.94%  '
$ws.Range('D36').Value = "'0.9997"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.93%  '
$ws.Range('D37').Value = "'0.9345"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.79%  '
$ws.Range('D38').Value = "'114.10"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +15.03%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = "'1.989"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.51%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = "'2.435"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.11%  '
$ws.Range('E41').Value = '  +0.39%  '
$ws.Range('D42').Value = "'5.791"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +18.39%  '
$ws.Range('E43').Value = '  +2.48%  '
$ws.Range('D44').Value = "'0.3922"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.78%  '
$ws.Range('D45').Value = "'0.1219"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +9.47%  '
$ws.Range('D46').Value = "'6.759"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +9.20%  '
$ws.Range('D47').Value = "'0.05337"
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Value = "'7.941"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.67%  '
$ws.Range('D49').Value = "'30.79"
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Value = "'1.264"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.05%  '
$ws.Range('D51').Value = "'0.3456"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.94%  '
